# Update the "想去人数" (interested-count) values in column F
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value  = 156
    $ws.Range("F3").Value  = 1722
    $ws.Range("F5").Value  = 1126
    $ws.Range("F7").Value  = 12012
    $ws.Range("F10").Value = 479
    $ws.Range("F12").Value = 1113
    $ws.Range("F13").Value = 867
    $ws.Range("F14").Value = 13486
    $ws.Range("F15").Value = 13492
    $ws.Range("F20").Value = 488
    $ws.Range("F21").Value = 96
    $ws.Range("F24").Value = 176
}

# F23 differs between the two sheets in both old and new value,
# so it needs to be set individually per sheet.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F23").Value = 1716

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F23").Value = 1718
